$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source file gains a header row's worth of real column titles, three
# new trailing columns (K, L, M), and one extra data row (old warehouse
# code 301 that was missing from the extract). Concretely: a new row is
# inserted right under the header so every existing data row shifts down
# by one, then the header row is rewritten as column titles, the newly
# freed row 2 is populated with the missing "301" record, and every data
# row gets its Warehouse group code / Extraction date / Replacement
# multiplier columns (G/H/J) filled in.

# Insert a new blank row at position 2; this pushes the old rows 2-21
# (the "566024A" / "566085A" replacement records) down to rows 3-22.
$ws.Rows(2).Insert()

# Rewrite row 1 as the header row.
$ws.Cells.Item(1,1).Value = "Row transaction type"
$ws.Cells.Item(1,2).Value = "Replacing Item code"
$ws.Cells.Item(1,3).Value = "Replacing prefix"
$ws.Cells.Item(1,4).Value = "Replaced Item code"
$ws.Cells.Item(1,5).Value = "Replaced prefix"
$ws.Cells.Item(1,6).Value = "Warehouse code"
$ws.Cells.Item(1,7).Value = "Warehouse group code"
$ws.Cells.Item(1,8).Value = "Extraction date"
$ws.Cells.Item(1,9).Value = "Inherit stock"
$ws.Cells.Item(1,10).Value = "Replacement multiplier"
$ws.Cells.Item(1,11).Value = "Replacement Description"
$ws.Cells.Item(1,12).Value = "Free text 1"
$ws.Cells.Item(1,13).Value = "Free text 2"

# Populate the newly inserted row 2 with the missing "301" warehouse
# record, matching the pattern shared by every other data row.
$ws.Cells.Item(2,1).Value = "M"
$ws.Cells.Item(2,2).Value = "566024A"
$ws.Cells.Item(2,4).Value = "566085A"
$ws.Cells.Item(2,6).Value = 301
$ws.Cells.Item(2,7).Value = "ERA"
$ws.Cells.Item(2,8).Value = "20220525T1548"
$ws.Cells.Item(2,10).Value = 1

# Fill in Warehouse group code (G), Extraction date (H) and Replacement
# multiplier (J) for all the shifted data rows (now rows 3-22).
for ($r = 3; $r -le 22; $r++) {
    $ws.Cells.Item($r,7).Value = "ERA"
    $ws.Cells.Item($r,8).Value = "20220525T1548"
    $ws.Cells.Item($r,10).Value = 1
}
